$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# --- Remove the obsolete VRT annotation shapes on slide 1 (superseded by the new flow chart) ---
$idsToDelete = @(303, 304, 305, 306, 308, 310, 311)
foreach ($id in $idsToDelete) {
    $sh = Get-ShapeById $s.Shapes $id
    if ($sh -ne $null) {
        $sh.Delete()
    }
}

# --- Reposition the "Rainfall [mm]" label (shape 301) that used to sit below the deleted shapes ---
$sh301 = Get-ShapeById $s.Shapes 301
if ($sh301 -ne $null) {
    $sh301.Top = 228.9055
}

# --- Refresh the date / slide-number footer fields on the slide master and every slide layout ---
function Update-DateSlideNum($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Type -eq 14 -and $sh.HasTextFrame) {
            $phType = $sh.PlaceholderFormat.Type
            if ($phType -eq 16) {
                $sh.TextFrame.TextRange.Text = "02/11/2023"
            } elseif ($phType -eq 13) {
                $sh.TextFrame.TextRange.Text = [char]0x2039 + "#" + [char]0x203A
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateSlideNum $master.Shapes

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DateSlideNum $layouts.Item($i).Shapes
}

Write-Host "Edit complete"
